# Add a new ON (obligación negociable) row for TLCPD (TELECOM) just above
# the existing TSC3D row, shifting every row from 78 down one position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 78 - Excel copies formatting from the row above,
# which is exactly the s="2"/"3" style pattern the target row needs.
$ws.Rows.Item(78).Insert()

# Populate the new row 78 with the new bond's data.
$ws.Cells.Item(78, 1).Value  = "TLCPD"                  # name
$ws.Cells.Item(78, 2).Value  = "TELECOM"                # empresa
$ws.Cells.Item(78, 3).Value  = "CCL"                    # curr
$ws.Cells.Item(78, 4).Value  = "NY"                      # law
$ws.Cells.Item(78, 5).Value  = 45805                     # start_date (28/5/2025)
$ws.Cells.Item(78, 6).Value  = 48727                     # end_date   (28/5/2033)
$ws.Cells.Item(78, 7).Value  = 6                         # payment_frequency
$ws.Cells.Item(78, 8).Value  = "28/5/2032;28/5/2033"     # amortization_dates
$ws.Cells.Item(78, 9).Value  = "50;50"                   # amortizations
$ws.Cells.Item(78, 10).Value = 0.095                     # rate
$ws.Cells.Item(78, 11).Value = 100                        # outstanding
$ws.Cells.Item(78, 12).Value = "AA+"                      # calificación

# Re-apply the AutoFilter over the grown range (A1:L103) - turning it off
# first, since toggling AutoFilter on an already-filtered sheet removes it.
$ws.AutoFilterMode = $false
$ws.Range("A1:L103").AutoFilter()

# The _FilterDatabase defined name needs to track the same grown range.
$wb.Names.Item("listado_ons!_FilterDatabase").RefersTo = "=listado_ons!`$A`$1:`$L`$103"

# Move the view: freeze pane scroll + selection on the bottom-right pane,
# matching where the author was working (row 77/78, column L).
$win = $excel.ActiveWindow
$win.Panes.Item($win.Panes.Count).ScrollRow = 52
$win.Panes.Item($win.Panes.Count).ScrollColumn = 5
$ws.Range("L77:L78").Select()
